# Retry implementation fetching data
# The underlying fetch/aggregation logic was retried, yielding row counts
# that are now 3x the previously recorded values in column C ("count").
# Multiply every numeric value in C2:C196 by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 196

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value = ($val * 3)
    }
}
